$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "47.11B"
$ws.Range("B10").Value = "Commerce d'alimentation générale"
$ws.Range("A10:B10").Style = $ws.Range("A9:B9").Style
